# Auto-generated PowerShell Excel COM-interop script
# Adds Corequisites / Concurrent / Recommended columns (D, E, F) and moves the
# "Terms Typically Offered" data into column G, splitting out embedded
# "Recommended:" / "Corequisite:" text that had been concatenated into the old
# column C (Prerequisites) / D (Terms Typically Offered) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=1; D="Corequisites"; E="Concurrent"; F="Recommended"; G="Terms Typically Offered" },
    @{ Row=2; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=3; D="NA"; E="NA"; F="NA"; G="F, SP" },
    @{ Row=4; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=5; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=6; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=7; C="MATH 141 with grade C- or better and MATH 142 or MATH 182 (or concurrent enrollment)."; D="NA"; E="NA"; F="High School Physics."; G="TBD " },
    @{ Row=8; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=9; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=10; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=11; C="Completion of GE Area A1 with a grade of C- or better."; D="NA"; E="NA"; F="Completion of GE Area A2."; G="F, W, SP " },
    @{ Row=12; C="Completion of GE Area A1 with a grade of C- or better."; D="NA"; E="NA"; F="Completion of GE Area A2."; G="TBD " },
    @{ Row=13; C="Completion of GE Area A1 with a grade of C- or better. For Engineering students only."; D="NA"; E="NA"; F="Completion of GE Area A2."; G="SP " },
    @{ Row=14; D="NA"; E="NA"; F="NA"; G="F" },
    @{ Row=15; D="NA"; E="NA"; F="NA"; G="W" },
    @{ Row=16; D="NA"; E="NA"; F="NA"; G="SP" },
    @{ Row=17; D="NA"; E="NA"; F="NA"; G="F, W, SP" },
    @{ Row=18; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=19; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=20; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=21; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=22; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=23; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=24; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=25; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=26; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=27; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=28; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=29; C="Completion of GE Area A with grades of C- or better; or for PHIL majors GE Area A3 with a grade of C- or better."; D="NA"; E="NA"; F="PHIL 126."; G="F " },
    @{ Row=30; C="Completion of GE Area A with grades of C- or better; or for PHIL majors GE Area A3 with a grade of C- or better."; D="NA"; E="NA"; F="PHIL 126."; G="W, SP " },
    @{ Row=31; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=32; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=33; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=34; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=35; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=36; D="NA"; E="NA"; F="Sophomore standing."; G="F, W, SP" },
    @{ Row=37; D="NA"; E="NA"; F="NA"; G="W" },
    @{ Row=38; D="NA"; E="NA"; F="NA"; G="SP" },
    @{ Row=39; D="NA"; E="NA"; F="NA"; G="W, SP" },
    @{ Row=40; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=41; D="NA"; E="NA"; F="NA"; G="F, W, SP" },
    @{ Row=42; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=43; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=44; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=45; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=46; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=47; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=48; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=49; C="Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1, C2, or C3."; D="NA"; E="NA"; F="Completion of GE Area C1 or C3."; G="TBD " },
    @{ Row=50; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=51; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=52; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=53; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=54; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=55; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=56; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=57; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=58; C="Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D1."; D="NA"; E="NA"; F="POLS 112 (GE Area D1)."; G="TBD " },
    @{ Row=59; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=60; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=61; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=62; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=63; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=64; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=65; C="Junior standing or History major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D1, D2, or D3."; D="NA"; E="NA"; F="One or more courses in GE Area B."; G="TBD " },
    @{ Row=66; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=67; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=68; D="NA"; E="NA"; F="HNRS 265."; G="W, SP" },
    @{ Row=69; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=70; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=71; C="Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4."; D="NA"; E="NA"; F="UNIV 391 and completion of GE Areas D2 and D3."; G="TBD " },
    @{ Row=72; D="NA"; E="NA"; F="NA"; G="F, W, SP" },
    @{ Row=73; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=74; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=75; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=76; C="one of the PHYS 104; PHYS 118; PHYS 121; or PHYS 141; and junior standing."; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=77; D="Major-specific senior project."; E="NA"; F="NA"; G="SP" },
    @{ Row=78; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=79; D="NA"; E="NA"; F="NA"; G="TBD" },
    @{ Row=80; C="Junior standing and completion of GE Area B, or graduate standing."; D="NA"; E="NA"; F="UNIV 391, GE Area D2, and GE Area D3."; G="TBD " },
    @{ Row=81; D="NA"; E="NA"; F="NA"; G="TBD" }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

Write-Output "Updated $($rows.Count) rows with Corequisites/Concurrent/Recommended/Terms columns."
